$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 105: correct the date/time value in column A ---
# (was 45488.6178819444, corrected to 45488.2916666667)
$ws.Range("A105").Value = 45488.2916666667

# --- Row 106: new row of results appended by the R script ---
$ws.Range("B106").Value = 300
$ws.Range("C106").Value = 6
$ws.Range("D106").Value = 6
$ws.Range("E106").Value = 6
$ws.Range("F106").Value = 6
$ws.Range("H106").Value = "PAL.MI"

# A106 needs the same date/time number format as the rest of column A.
# Setting .Value alone would leave it on the default "General" style, so
# copy the formatting down from the cell above (reuses the existing date
# style rather than registering a new one).
$ws.Range("A106").Value = 45489.560787037
$ws.Range("A105").Copy() | Out-Null
$ws.Range("A106").PasteSpecial(-4122) | Out-Null

# G106 holds "6" as text (matching the rest of column G, which stores
# adj_close as text). Entering it directly would auto-convert to a number,
# so build it with a formula that forces text, then flatten it down to a
# plain value via copy/paste so no formula or extra style is left behind.
$ws.Range("G106").Formula = '="6"'
$ws.Range("G106").Copy() | Out-Null
$ws.Range("G106").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0
